$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 19.539082
$ws.Range("H2").Value = 58.61724600000001
$ws.Range("I2").Value = 0.224220971665117
$ws.Range("J2").Value = 0.224220971665117
$ws.Range("M2").Value = 61.156892
$ws.Range("N2").Value = 183.470676
$ws.Range("O2").Value = 0.9308124486389074
$ws.Range("P2").Value = 0.9308124486389074
$ws.Range("Q2").Value = 1194.949527653144
$ws.Range("R2").Value = 10754.5457488783
$ws.Range("S2").Value = 0.2087076716718026
$ws.Range("T2").Value = 0.2087076716718026
$ws.Range("G3").Value = 19.539082
$ws.Range("H3").Value = 58.61724600000001
$ws.Range("I3").Value = 0.224220971665117
$ws.Range("J3").Value = 0.224220971665117
$ws.Range("O3").Value = 0.02171808228502914
$ws.Range("P3").Value = 0.02171808228502914
$ws.Range("Q3").Value = 27.88103253880667
$ws.Range("R3").Value = 250.92929284926
$ws.Range("S3").Value = 0.004869649512652198
$ws.Range("T3").Value = 0.004869649512652197
$ws.Range("G4").Value = 19.539082
$ws.Range("H4").Value = 58.61724600000001
$ws.Range("I4").Value = 0.224220971665117
$ws.Range("J4").Value = 0.224220971665117
$ws.Range("M4").Value = 2.00294
$ws.Range("N4").Value = 6.00882
$ws.Range("O4").Value = 0.03048489589491914
$ws.Range("P4").Value = 0.03048489589491914
$ws.Range("Q4").Value = 39.13560890108001
$ws.Range("R4").Value = 352.22048010972
$ws.Range("S4").Value = 0.006835352978668705
$ws.Range("T4").Value = 0.006835352978668705
$ws.Range("G5").Value = 19.539082
$ws.Range("H5").Value = 58.61724600000001
$ws.Range("I5").Value = 0.224220971665117
$ws.Range("J5").Value = 0.224220971665117
$ws.Range("M5").Value = 1.115932333333334
$ws.Range("N5").Value = 3.347797
$ws.Range("O5").Value = 0.01698457318114416
$ws.Range("P5").Value = 0.01698457318114415
$ws.Range("Q5").Value = 21.80429336745134
$ws.Range("R5").Value = 196.238640307062
$ws.Range("S5").Value = 0.003808297501993429
$ws.Range("T5").Value = 0.003808297501993429
$ws.Range("I6").Value = 0.3010605798326856
$ws.Range("J6").Value = 0.3010605798326856
$ws.Range("M6").Value = 61.156892
$ws.Range("N6").Value = 183.470676
$ws.Range("O6").Value = 0.9308124486389074
$ws.Range("P6").Value = 0.9308124486389074
$ws.Range("Q6").Value = 1604.45383406577
$ws.Range("R6").Value = 14440.08450659193
$ws.Range("S6").Value = 0.2802309355027114
$ws.Range("T6").Value = 0.2802309355027114
$ws.Range("I7").Value = 0.3010605798326856
$ws.Range("J7").Value = 0.3010605798326856
$ws.Range("O7").Value = 0.02171808228502914
$ws.Range("P7").Value = 0.02171808228502914
$ws.Range("S7").Value = 0.006538458445584852
$ws.Range("T7").Value = 0.006538458445584851
$ws.Range("I8").Value = 0.3010605798326856
$ws.Range("J8").Value = 0.3010605798326856
$ws.Range("M8").Value = 2.00294
$ws.Range("N8").Value = 6.00882
$ws.Range("O8").Value = 0.03048489589491914
$ws.Range("P8").Value = 0.03048489589491914
$ws.Range("Q8").Value = 52.54722169994667
$ws.Range("R8").Value = 472.92499529952
$ws.Range("S8").Value = 0.009177800434263415
$ws.Range("T8").Value = 0.009177800434263415
$ws.Range("I9").Value = 0.3010605798326856
$ws.Range("J9").Value = 0.3010605798326856
$ws.Range("M9").Value = 1.115932333333334
$ws.Range("N9").Value = 3.347797
$ws.Range("O9").Value = 0.01698457318114416
$ws.Range("P9").Value = 0.01698457318114415
$ws.Range("Q9").Value = 29.27653535393245
$ws.Range("R9").Value = 263.488818185392
$ws.Range("S9").Value = 0.005113385450125942
$ws.Range("T9").Value = 0.005113385450125942
$ws.Range("G10").Value = 19.67155566666667
$ws.Range("H10").Value = 59.014667
$ws.Range("I10").Value = 0.2257411748281949
$ws.Range("J10").Value = 0.2257411748281949
$ws.Range("M10").Value = 61.156892
$ws.Range("N10").Value = 183.470676
$ws.Range("O10").Value = 0.9308124486389074
$ws.Range("P10").Value = 0.9308124486389074
$ws.Range("Q10").Value = 1203.051205378321
$ws.Range("R10").Value = 10827.46084840489
$ws.Range("S10").Value = 0.2101226957004558
$ws.Range("T10").Value = 0.2101226957004558
$ws.Range("G11").Value = 19.67155566666667
$ws.Range("H11").Value = 59.014667
$ws.Range("I11").Value = 0.2257411748281949
$ws.Range("J11").Value = 0.2257411748281949
$ws.Range("O11").Value = 0.02171808228502914
$ws.Range("P11").Value = 0.02171808228502914
$ws.Range("Q11").Value = 28.07006407114111
$ws.Range("R11").Value = 252.63057664027
$ws.Range("S11").Value = 0.004902665410037887
$ws.Range("T11").Value = 0.004902665410037886
$ws.Range("G12").Value = 19.67155566666667
$ws.Range("H12").Value = 59.014667
$ws.Range("I12").Value = 0.2257411748281949
$ws.Range("J12").Value = 0.2257411748281949
$ws.Range("M12").Value = 2.00294
$ws.Range("N12").Value = 6.00882
$ws.Range("O12").Value = 0.03048489589491914
$ws.Range("P12").Value = 0.03048489589491914
$ws.Range("Q12").Value = 39.40094570699333
$ws.Range("R12").Value = 354.60851136294
$ws.Range("S12").Value = 0.006881696213834264
$ws.Range("T12").Value = 0.006881696213834263
$ws.Range("G13").Value = 19.67155566666667
$ws.Range("H13").Value = 59.014667
$ws.Range("I13").Value = 0.2257411748281949
$ws.Range("J13").Value = 0.2257411748281949
$ws.Range("M13").Value = 1.115932333333334
$ws.Range("N13").Value = 3.347797
$ws.Range("O13").Value = 0.01698457318114416
$ws.Range("P13").Value = 0.01698457318114415
$ws.Range("Q13").Value = 21.95212501539989
$ws.Range("R13").Value = 197.569125138599
$ws.Range("S13").Value = 0.003834117503866934
$ws.Range("T13").Value = 0.003834117503866933
$ws.Range("G14").Value = 21.69639766666667
$ws.Range("H14").Value = 65.08919299999999
$ws.Range("I14").Value = 0.2489772736740025
$ws.Range("J14").Value = 0.2489772736740025
$ws.Range("M14").Value = 61.156892
$ws.Range("N14").Value = 183.470676
$ws.Range("O14").Value = 0.9308124486389074
$ws.Range("P14").Value = 0.9308124486389074
$ws.Range("Q14").Value = 1326.884248889385
$ws.Range("R14").Value = 11941.95824000447
$ws.Range("S14").Value = 0.2317511457639376
$ws.Range("T14").Value = 0.2317511457639376
$ws.Range("G15").Value = 21.69639766666667
$ws.Range("H15").Value = 65.08919299999999
$ws.Range("I15").Value = 0.2489772736740025
$ws.Range("J15").Value = 0.2489772736740025
$ws.Range("O15").Value = 0.02171808228502914
$ws.Range("P15").Value = 0.02171808228502914
$ws.Range("Q15").Value = 30.95938536514777
$ws.Range("R15").Value = 278.6344682863299
$ws.Range("S15").Value = 0.005407308916754206
$ws.Range("T15").Value = 0.005407308916754205
$ws.Range("G16").Value = 21.69639766666667
$ws.Range("H16").Value = 65.08919299999999
$ws.Range("I16").Value = 0.2489772736740025
$ws.Range("J16").Value = 0.2489772736740025
$ws.Range("M16").Value = 2.00294
$ws.Range("N16").Value = 6.00882
$ws.Range("O16").Value = 0.03048489589491914
$ws.Range("P16").Value = 0.03048489589491914
$ws.Range("Q16").Value = 43.45658274247334
$ws.Range("R16").Value = 391.10924468226
$ws.Range("S16").Value = 0.007590046268152758
$ws.Range("T16").Value = 0.007590046268152757
$ws.Range("G17").Value = 21.69639766666667
$ws.Range("H17").Value = 65.08919299999999
$ws.Range("I17").Value = 0.2489772736740025
$ws.Range("J17").Value = 0.2489772736740025
$ws.Range("M17").Value = 1.115932333333334
$ws.Range("N17").Value = 3.347797
$ws.Range("O17").Value = 0.01698457318114416
$ws.Range("P17").Value = 0.01698457318114415
$ws.Range("Q17").Value = 24.21171167309123
$ws.Range("R17").Value = 217.905405057821
$ws.Range("S17").Value = 0.004228772725157852
$ws.Range("T17").Value = 0.004228772725157851
